$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3176.7693
$ws.Range("I29").Value = 2828.8572
$ws.Range("J29").Value = 3582.6667
$ws.Range("K29").Value = 8486.571599999999
$ws.Range("L29").Value = 10748.0001
$ws.Range("M29").Value = -8205.571599999999
$ws.Range("N29").Value = -11310.0001

$ws.Range("H32").Value = 4583.3687
$ws.Range("I32").Value = 4096.5713
$ws.Range("J32").Value = 4867.3335
$ws.Range("K32").Value = 4096.5713
$ws.Range("L32").Value = 4867.3335
$ws.Range("M32").Value = -3770.5713
$ws.Range("N32").Value = -5519.3335

$ws.Range("H98").Value = 1051.421
$ws.Range("I98").Value = 980.64703
$ws.Range("K98").Value = 980.64703
$ws.Range("M98").Value = 517.35297

$ws.Range("H106").Value = 3640.5715
$ws.Range("I106").Value = 3415
$ws.Range("K106").Value = 3415
$ws.Range("M106").Value = -2784

$ws.Range("H122").Value = 1051.421
$ws.Range("I122").Value = 980.64703
$ws.Range("K122").Value = 2941.94109
$ws.Range("M122").Value = -491.9410899999998

$ws.Range("H132").Value = 2492.9656
$ws.Range("I132").Value = 2138.7693
$ws.Range("K132").Value = 6416.3079
$ws.Range("M132").Value = -3886.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11369633
$ws.Range("I32").Value = 15875325
$ws.Range("J32").Value = 15290.2
$ws.Range("K32").Value = 15875325
$ws.Range("L32").Value = 15290.2
$ws.Range("M32").Value = -15875038
$ws.Range("N32").Value = -15864.2

$ws.Range("H45").Value = 1648.8
$ws.Range("I45").Value = 1665.6666
$ws.Range("K45").Value = 1665.6666
$ws.Range("M45").Value = -1288.6666

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H110").Value = 12341.97
$ws.Range("I110").Value = 12903.333
$ws.Range("J110").Value = 6728.3335
$ws.Range("K110").Value = 12903.333
$ws.Range("L110").Value = 6728.3335
$ws.Range("M110").Value = -10858.333
$ws.Range("N110").Value = -10818.3335

$ws.Range("H132").Value = 26318332
$ws.Range("I132").Value = 1907.9678
$ws.Range("K132").Value = 5723.903399999999
$ws.Range("M132").Value = -3193.903399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18553.525
$ws.Range("I86").Value = 12897.223
$ws.Range("K86").Value = 12897.223
$ws.Range("M86").Value = -11774.223

$ws.Range("H89").Value = 18553.525
$ws.Range("I89").Value = 12897.223
$ws.Range("K89").Value = 64486.115
$ws.Range("M89").Value = -58870.115

$ws.Range("H105").Value = 9353.706
$ws.Range("I105").Value = 15315.625
$ws.Range("J105").Value = 4054.2222
$ws.Range("K105").Value = 15315.625
$ws.Range("L105").Value = 4054.2222
$ws.Range("M105").Value = -13568.625
$ws.Range("N105").Value = -7548.2222

$ws.Range("H107").Value = 3692.3333
$ws.Range("I107").Value = 3528.875
$ws.Range("K107").Value = 3528.875
$ws.Range("M107").Value = -1608.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 20000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 20000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 20000
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -20280

$ws.Range("H31").Value = 24394734
$ws.Range("I31").Value = 2944.6177
$ws.Range("J31").Value = 142869140
$ws.Range("K31").Value = 2944.6177
$ws.Range("L31").Value = 142869140
$ws.Range("M31").Value = -2649.6177
$ws.Range("N31").Value = -142869730

$ws.Range("H34").Value = 24394734
$ws.Range("I34").Value = 2944.6177
$ws.Range("J34").Value = 142869140
$ws.Range("K34").Value = 2944.6177
$ws.Range("L34").Value = 142869140
$ws.Range("M34").Value = -2742.6177
$ws.Range("N34").Value = -142869544

$ws.Range("H58").Value = 2503.5833
$ws.Range("I58").Value = 1930.8422
$ws.Range("J58").Value = 4680
$ws.Range("K58").Value = 1930.8422
$ws.Range("L58").Value = 4680
$ws.Range("M58").Value = -1727.8422
$ws.Range("N58").Value = -5086

$ws.Range("H99").Value = 10329.052
$ws.Range("J99").Value = 11219.218
$ws.Range("L99").Value = 11219.218
$ws.Range("N99").Value = -14215.218

$ws.Range("H126").Value = 10329.052
$ws.Range("J126").Value = 11219.218
$ws.Range("L126").Value = 33657.654
$ws.Range("N126").Value = -38597.654

$ws.Range("H132").Value = 3530.8845
$ws.Range("J132").Value = 6204.6665
$ws.Range("L132").Value = 18613.9995
$ws.Range("N132").Value = -23673.9995

$ws.Range("H135").Value = 76750
$ws.Range("J135").Value = 76750
$ws.Range("L135").Value = 76750
$ws.Range("N135").Value = -86890

$ws.Range("H136").Value = 2503.5833
$ws.Range("I136").Value = 1930.8422
$ws.Range("J136").Value = 4680
$ws.Range("K136").Value = 5792.5266
$ws.Range("L136").Value = 14040
$ws.Range("M136").Value = -3242.5266
$ws.Range("N136").Value = -19140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 16949
$ws.Range("I14").Value = 16949
$ws.Range("K14").Value = 50847
$ws.Range("M14").Value = -50674

$ws.Range("H69").Value = 4151.25
$ws.Range("J69").Value = 4999.8335
$ws.Range("L69").Value = 14999.5005
$ws.Range("N69").Value = -16621.5005

$ws.Range("H72").Value = 4151.25
$ws.Range("J72").Value = 4999.8335
$ws.Range("L72").Value = 44998.5015
$ws.Range("N72").Value = -53110.5015

$ws.Range("H140").Value = 1923.1666
$ws.Range("I140").Value = 1701.2
$ws.Range("K140").Value = 5103.6
$ws.Range("M140").Value = 76.39999999999964

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1743.0714
$ws.Range("I80").Value = 1573
$ws.Range("K80").Value = 1573
$ws.Range("M80").Value = -575

$ws.Range("H83").Value = 1743.0714
$ws.Range("I83").Value = 1573
$ws.Range("K83").Value = 7865
$ws.Range("M83").Value = -2873

$ws.Range("H93").Value = 36247.25
$ws.Range("I93").Value = 32000
$ws.Range("J93").Value = 40494.5
$ws.Range("K93").Value = 32000
$ws.Range("L93").Value = 40494.5
$ws.Range("M93").Value = -30128
$ws.Range("N93").Value = -44238.5

$ws.Range("H107").Value = 421.95
$ws.Range("I107").Value = 766.125
$ws.Range("J107").Value = 192.5
$ws.Range("K107").Value = 766.125
$ws.Range("L107").Value = 192.5
$ws.Range("M107").Value = 1153.875
$ws.Range("N107").Value = -4032.5

$ws.Range("H122").Value = 1860.421
$ws.Range("I122").Value = 1459.875
$ws.Range("K122").Value = 4379.625
$ws.Range("M122").Value = -1929.625

$ws.Range("H132").Value = 2979
$ws.Range("I132").Value = 3058.9167
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 9176.750100000001
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -6646.750100000001
$ws.Range("N132").Value = -12558.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1498.75
$ws.Range("I13").Value = 1498.3334
$ws.Range("J13").Value = 1500
$ws.Range("K13").Value = 1498.3334
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = -1358.3334
$ws.Range("N13").Value = -1780

$ws.Range("H44").Value = 17873.75
$ws.Range("I44").Value = 15750
$ws.Range("J44").Value = 19997.5
$ws.Range("K44").Value = 15750
$ws.Range("L44").Value = 19997.5
$ws.Range("M44").Value = -15196
$ws.Range("N44").Value = -21105.5

$ws.Range("H62").Value = 7637.6924
$ws.Range("J62").Value = 7662.727
$ws.Range("L62").Value = 7662.727
$ws.Range("N62").Value = -8910.726999999999

$ws.Range("H65").Value = 7637.6924
$ws.Range("J65").Value = 7662.727
$ws.Range("L65").Value = 38313.635
$ws.Range("N65").Value = -44553.635

$ws.Range("H122").Value = 47668508
$ws.Range("I122").Value = 58883576
$ws.Range("J122").Value = 4462.25
$ws.Range("K122").Value = 176650728
$ws.Range("L122").Value = 13386.75
$ws.Range("M122").Value = -176648278
$ws.Range("N122").Value = -18286.75
